# JORequestForm.xlsx content update
# Updates the request header fields, the general-description line, and the
# first job-order line item; clears the now-unused sample line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header block (rows 7-12) ---
# "JO Request to:" value
$ws.Range("C7").Value = "bacolod"

# Dates (Date Prepared / Completion Date / Delivery Date) -> 2022-10-07
$ws.Range("C8").Value = 44841
$ws.Range("I8").Value = 44841
$ws.Range("I9").Value = 44841

# "Department:" value (now points users to the code sheet)
$ws.Range("C9").Value = "REFER TO DEPARTMENT CODE SHEET"

# "JO No.:" value
$ws.Range("C10").Value = "ADM 1001"

# "Requested By:" value
$ws.Range("C11").Value = "Stephine"

# "Purpose:" value
$ws.Range("C12").Value = "test"

# --- General description banner (row 14) ---
$ws.Range("A14").Value = "GENERAL DESCRIPTION"

# --- Line item 1 (row 15): scope of work / qty / uom ---
$ws.Range("B15").Value = " Labor`nA. Replace advise parts`na.1. Rack End`na.2. Tie Rod End`nB. Wheel Alignment"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = ""

# --- Clear the now-unused sample line items (rows 16 & 17) ---
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""

$ws.Range("A17").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""
